$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.180.01"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.813.15"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'312.14"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.4616"
$ws.Range("E7").Value = "  +4.76%  "
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "'0.07387"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "'0.8646"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "'20.59"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.813.56"
$ws.Range("E12").Value = "  -6.45%  "
$ws.Range("D13").Value = "'6.649"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "'5.383"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "'0.07081"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "'91.73"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'0.000008730"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'14.86"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").Value = "27.179.95"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "'5.304"
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("D23").Value = "'10.89"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "2.042.81"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").Value = "'1.928"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "'151.84"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'2.225"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").Value = "'18.47"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "'5.264"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "'116.76"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'0.08885"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'0.7710"
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("D33").Value = "'1.171"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "'4.515"
$ws.Range("E34").Value = "  +2.26%  "
$ws.Range("D35").Value = "'2.922"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'1.113"
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("D38").Value = "'0.01959"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'0.05231"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("D40").Value = "'7.244"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").Value = "'2.914"
$ws.Range("E41").Value = "  +3.96%  "
$ws.Range("D42").Value = "'2.370"
$ws.Range("E42").Value = "  +20.68%  "
$ws.Range("D43").Value = "'0.5275"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'8.567"
$ws.Range("D46").Value = "'0.5014"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'10.40"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'105.18"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'1.667"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "'0.06320"
$ws.Range("E51").Value = "  +0.24%  "
